$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text such as "70.119.52" or "1.00".
# Several of the refreshed prices look like ordinary decimals (e.g.
# "604.30", "1.00"); writing them with a bare .Value assignment would let
# Excel reinterpret them as numbers and mangle the text (dropping a
# trailing zero, turning "1.00" into 1, etc). Force those specific cells
# to a text number format first so the original text formatting survives.

$ws.Range("D2").Value = '70.119.52'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '3.542.89'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.09'
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '197.08'
$ws.Range("E6").Value = '  +5.87%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.209'
$ws.Range("E9").Value = '  -3.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.654'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.06'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '4.110.05'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '604.30'
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.25'
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '70.209.18'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = '3.548.68'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.996'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.02'
$ws.Range("E22").Value = '  +2.97%  '
$ws.Range("E23").Value = '  +5.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.78'
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.15'
$ws.Range("E26").Value = '  +3.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.02'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.65'
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.80'
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  +20.91%  '
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.66'
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").Value = '0.0₃0843'
$ws.Range("E35").Value = '  +8.28%  '
$ws.Range("D36").Value = '3.779.57'
$ws.Range("E36").Value = '  +6.97%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.08'
$ws.Range("E37").Value = '  -3.91%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.68'
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.81'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '490.23'
$ws.Range("E42").Value = '  -7.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.134'
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0458'
$ws.Range("E44").Value = '  -1.48%  '
$ws.Range("E45").Value = '  -3.59%  '
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.29'
$ws.Range("E47").Value = '  -2.30%  '
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("E49").Value = '  -4.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000250'
$ws.Range("E50").Value = '  +3.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.45'
$ws.Range("E51").Value = '  -2.28%  '
